$d = $word.ActiveDocument

# --- Part 1: add the new sentence right after "(ref L. Julia)" -------------
# Locate the existing "(ref L. Julia)" text and collapse the found range to
# its end point (immediately before the hidden _GoBack bookmark).
$rng = $d.Content
$found = $rng.Find.Execute("(ref L. Julia)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)

$rng.InsertAfter(".")
$rng.Collapse(0)
$rng.InsertAfter(" Si cette fois, nous comparons le terme ")
$rng.Collapse(0)
$rng.InsertAfter("« ")
$rng.Collapse(0)
$rng.InsertAfter("intelligence artificielle forte")
$rng.Collapse(0)
$rng.InsertAfter(" »")
$rng.Collapse(0)
$rng.InsertAfter(" sur me moteur de recherche google, nous obtenons XXX entre 2004 et XXX ")
$rng.Collapse(0)
$rng.InsertAfter("à")
$rng.Collapse(0)
$rng.InsertAfter(" aujourd")
$rng.Collapse(0)
$rng.InsertAfter("’")
$rng.Collapse(0)
$rng.InsertAfter("hui.")
$rng.Collapse(0)

# The paragraph used to end with the bookmark followed by a separate run
# containing only ".". Now that our new text supplies its own closing
# period, that old stray "." run (sitting right after the _GoBack bookmark)
# must be removed so the paragraph doesn't end in "..".
$rng.MoveEnd(1, 2)
$rng.Delete()

# --- Part 2: lastRenderedPageBreak shifts from "sounds" to "image" ---------
# Adding the text above pushes the later content down, so the page break
# that used to fall right before "Explain the sounds in the database" now
# falls right before "Explain the image in the database" instead.
$rngImg = $d.Paragraphs.Item(29).Range
$rngImg.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' w:rsidR='00CC742C' w:rsidRDefault='00CC742C'><w:pPr><w:rPr><w:lang w:val='en-GB'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='en-GB'/></w:rPr><w:lastRenderedPageBreak/><w:t>Explain the image in the database</w:t></w:r></w:p>")

$rngSnd = $d.Paragraphs.Item(31).Range
$rngSnd.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' w:rsidR='00CC742C' w:rsidRDefault='00CC742C'><w:pPr><w:rPr><w:lang w:val='en-GB'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='en-GB'/></w:rPr><w:t>Explain the sounds in the database</w:t></w:r></w:p>")

Write-Host "All edits applied"
